$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 37
$ws.Range("D37").Value = "[Paper Review]Metric Learning for Adversarial Robustness"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1446&mod=document&pageid=1"

# Row 39
$ws.Range("D39").Value = "Anomaly Detection using Autoencoders"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Anomaly-Detection-using-Autoencoders-1"

# Row 45
$ws.Range("D45").Value = "데이터 분석 flow 에서의 유의해야 할 point, 주의할 점"
$ws.Range("E45").Value = "https://dive-into-ds.tistory.com/82"

# Row 51
$ws.Range("D51").Value = "[세이버메트릭스] 팀 득실점 비율과 승률 사이의 상관관계, 그리고 피타고리언 승률"
$ws.Range("E51").Value = "https://bskyvision.com/1131"
